$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are plain text in the source data (e.g. "1.00",
# European-grouped "51.104.72", subscript notation, etc). Excel's COM layer
# auto-converts a numeric-looking string assigned to .Value into a real
# number (dropping trailing zeros / introducing float noise), so for each
# D-column write we snapshot the cell style, force a Text number format,
# write the literal string, then restore the original style. This keeps the
# text verbatim while leaving the cell formatting untouched.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "51.104.72"
$ws.Range("E2").Value = "  -3.07%  "
Set-TextValue $ws.Range("D3") "2.903.96"
$ws.Range("E3").Value = "  -3.00%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.17%  "
Set-TextValue $ws.Range("D5") "370.23"
$ws.Range("E5").Value = "  +3.03%  "
Set-TextValue $ws.Range("D6") "102.85"
$ws.Range("E6").Value = "  -6.97%  "
$ws.Range("E7").Value = "  -5.80%  "
Set-TextValue $ws.Range("D8") "1.00"
$ws.Range("E9").Value = "  -7.26%  "
Set-TextValue $ws.Range("D10") "36.87"
$ws.Range("E10").Value = "  -6.68%  "
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("E12").Value = "  -5.38%  "
Set-TextValue $ws.Range("D13") "18.44"
$ws.Range("E13").Value = "  -5.73%  "
Set-TextValue $ws.Range("D14") "3.360.22"
$ws.Range("E14").Value = "  -3.08%  "
Set-TextValue $ws.Range("D15") "7.34"
$ws.Range("E15").Value = "  -6.51%  "
Set-TextValue $ws.Range("D16") "2.897.44"
$ws.Range("E16").Value = "  -4.86%  "
$ws.Range("E17").Value = "  -6.21%  "
Set-TextValue $ws.Range("D18") "51.035.07"
$ws.Range("E18").Value = "  -3.33%  "
$ws.Range("E19").Value = "  -6.89%  "
Set-TextValue $ws.Range("D20") "7.25"
$ws.Range("E20").Value = "  -5.68%  "
$ws.Range("E21").Value = "  -8.61%  "
Set-TextValue $ws.Range("D22") "0.0₃0944"
$ws.Range("E22").Value = "  -4.82%  "
Set-TextValue $ws.Range("D23") "68.32"
$ws.Range("E23").Value = "  -3.71%  "
Set-TextValue $ws.Range("D24") "259.70"
$ws.Range("E24").Value = "  -4.51%  "
Set-TextValue $ws.Range("D25") "2.70"
$ws.Range("E25").Value = "  -4.34%  "
Set-TextValue $ws.Range("D26") "0.168"
$ws.Range("E26").Value = "  -7.22%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -6.38%  "
Set-TextValue $ws.Range("D29") "7.10"
$ws.Range("E29").Value = "  -7.20%  "
Set-TextValue $ws.Range("D30") "0.103"
$ws.Range("E30").Value = "  -5.17%  "
$ws.Range("E31").Value = "  -6.37%  "
Set-TextValue $ws.Range("D32") "6.01"
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("E33").Value = "  -2.66%  "
Set-TextValue $ws.Range("D34") "34.94"
$ws.Range("E34").Value = "  -8.36%  "
Set-TextValue $ws.Range("D35") "51.16"
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("E36").Value = "  +0.34%  "
Set-TextValue $ws.Range("D37") "0.0419"
$ws.Range("E37").Value = "  -6.53%  "
Set-TextValue $ws.Range("D38") "3.05"
$ws.Range("E38").Value = "  -7.06%  "
Set-TextValue $ws.Range("D39") "2.64"
$ws.Range("E39").Value = "  -5.19%  "
Set-TextValue $ws.Range("D40") "16.96"
$ws.Range("E40").Value = "  -7.22%  "
Set-TextValue $ws.Range("D41") "1.85"
$ws.Range("E41").Value = "  -10.95%  "
$ws.Range("E42").Value = "  -6.36%  "
Set-TextValue $ws.Range("D43") "22.08"
$ws.Range("E43").Value = "  -7.97%  "
Set-TextValue $ws.Range("D44") "117.57"
$ws.Range("E44").Value = "  -1.95%  "
Set-TextValue $ws.Range("D45") "2.09"
$ws.Range("E45").Value = "  -3.81%  "
Set-TextValue $ws.Range("D46") "2.040.86"
$ws.Range("E46").Value = "  -5.45%  "
$ws.Range("E47").Value = "  -6.17%  "
$ws.Range("E48").Value = "  -9.70%  "
Set-TextValue $ws.Range("D49") "3.188.21"
$ws.Range("E49").Value = "  -3.03%  "
Set-TextValue $ws.Range("D50") "0.238"
$ws.Range("E50").Value = "  -3.67%  "
Set-TextValue $ws.Range("D51") "0.0311"
$ws.Range("E51").Value = "  -12.63%  "
